$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 655, shifting existing rows 655:696 down to 656:697
$ws.Rows.Item(655).Insert()

# Populate the newly inserted row with the new data point
# (leading apostrophe forces these date-look-alike / text values to stay
# literal text instead of being auto-converted to a date serial number;
# resetting the Style afterward drops the quote-prefix formatting so the
# cell ends up with the same plain/default style as its neighbours)
$ws.Range("A655").Value = "'2026/01/15"
$ws.Range("A655").Style = "Normal"
$ws.Range("B655").Value = "木"
$ws.Range("C655").Value = 13
$ws.Range("D655").Value = 201
